# Adds a new week of data (2022-01-17, serial 44578) for Melón (Calameño / Tuna)
# at Mercado Mayorista Lo Valledor de Santiago. The new observations are
# inserted as 9 new rows right before the existing row 902, which pushes all
# the rows that used to start at 902 down by 9 (902->911 ... 943->952).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows above the current row 902 (shifts old 902..943 -> 911..952)
$ws.Range("A902:A910").EntireRow.Insert()

# Common, constant values for every new row
$commonA = 6
$commonB = "Mercado Mayorista Lo Valledor de Santiago"
$commonC = "Metropolitana"
$commonDate = 44578
$commonE = 13
$commonF = 100112027
$commonG = "Melón"
$commonN = "`$/unidad"
$commonQ = 1
$commonR = "Hortaliza"

function Set-MelonRow($RowNum, $Variedad, $Calidad, $Mayor, $Menor, $MenorMax, $Promedio, $Region) {
    $ws.Cells.Item($RowNum, 1).Value = $commonA
    $ws.Cells.Item($RowNum, 2).Value = $commonB
    $ws.Cells.Item($RowNum, 3).Value = $commonC
    $ws.Cells.Item($RowNum, 4).Value = $commonDate
    $ws.Cells.Item($RowNum, 5).Value = $commonE
    $ws.Cells.Item($RowNum, 6).Value = $commonF
    $ws.Cells.Item($RowNum, 7).Value = $commonG
    $ws.Cells.Item($RowNum, 8).Value = $Variedad
    $ws.Cells.Item($RowNum, 9).Value = $Calidad
    $ws.Cells.Item($RowNum, 10).Value = $Mayor
    $ws.Cells.Item($RowNum, 11).Value = $Menor
    $ws.Cells.Item($RowNum, 12).Value = $MenorMax
    $ws.Cells.Item($RowNum, 13).Value = $Promedio
    $ws.Cells.Item($RowNum, 14).Value = $commonN
    $ws.Cells.Item($RowNum, 15).Value = $Region
    $ws.Cells.Item($RowNum, 16).Value = $Promedio
    $ws.Cells.Item($RowNum, 17).Value = $commonQ
    $ws.Cells.Item($RowNum, 18).Value = $commonR
}

Set-MelonRow 902 "Calameño" "Extra"   3500 750  800  771  "Región Metropolitana"
Set-MelonRow 903 "Calameño" "Primera" 4300 600  700  640  "Región Metropolitana"
Set-MelonRow 904 "Calameño" "Segunda" 3600 400  500  439  "Región Metropolitana"
Set-MelonRow 905 "Calameño" "Super"   3100 1000 1100 1042 "Región Metropolitana"
Set-MelonRow 906 "Tuna"     "Extra"   5100 750  800  772  "Región de O'Higgins"
Set-MelonRow 907 "Tuna"     "Primera" 7900 600  700  644  "Región de O'Higgins"
Set-MelonRow 908 "Tuna"     "Segunda" 7500 400  500  439  "Región de O'Higgins"
Set-MelonRow 909 "Tuna"     "Super"   4600 900  1000 941  "Región de O'Higgins"
Set-MelonRow 910 "Tuna"     "Tercera" 5900 200  300  244  "Región de O'Higgins"

Write-Host "done"
